$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 43053
$ws.Range("B7").Value = 0.31597222222222221
$ws.Range("C7").Value = 0.3576388888888889

$ws.Range("B7").NumberFormat = $ws.Range("B6").NumberFormat
$ws.Range("C7").NumberFormat = $ws.Range("C6").NumberFormat

$ws.Range("C8").Select()
